$wb = $excel.ActiveWorkbook

# --- Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A44").Value = 43935
$ws1.Range("A44").NumberFormat = $ws1.Range("A43").NumberFormat
$ws1.Range("B44").Value = 43
$ws1.Range("C44").Value = 124
$ws1.Range("D44").Value = 62
$ws1.Range("E44").Value = 176
$ws1.Range("F44").Value = 13
$ws1.Range("G44").Value = 66
$ws1.Range("H44").Value = 299
$ws1.Range("I44").Value = 4334
$ws1.Range("J44").Value = 54
$ws1.Range("K44").Value = 152
$ws1.Range("L44").Value = 634
$ws1.Range("M44").Value = 542
$ws1.Range("N44").Value = 826
$ws1.Range("O44").Value = 150
$ws1.Range("P44").Value = 390
$ws1.Range("Q44").Value = 7
$ws1.Range("R44").Value = 444
$ws1.Range("S44").Value = 8273

# --- Hoja2 ---
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("A44").Value = 43935
$ws2.Range("A44").NumberFormat = $ws2.Range("A43").NumberFormat
$ws2.Range("B44").Value = 43
$ws2.Range("C44").Value = 1
$ws2.Range("D44").Value = 0
$ws2.Range("E44").Value = 1
$ws2.Range("F44").Value = 0
$ws2.Range("G44").Value = 0
$ws2.Range("H44").Value = 2
$ws2.Range("I44").Value = 41
$ws2.Range("J44").Value = 0
$ws2.Range("K44").Value = 4
$ws2.Range("L44").Value = 8
$ws2.Range("M44").Value = 2
$ws2.Range("N44").Value = 21
$ws2.Range("O44").Value = 3
$ws2.Range("P44").Value = 5
$ws2.Range("Q44").Value = 0
$ws2.Range("R44").Value = 6
$ws2.Range("S44").Value = 94

# --- View / selection state (best effort mirror of the saved sheetViews) ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 34
$ws1.Range("A51:AC57").Select()
$ws1.Range("AC57").Activate()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 31
$ws2.Range("A51:AC57").Select()
$ws2.Range("A42").Select()
$ws1.Activate()
